$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Program")

# Update values for the Program worksheet (renamed test data entries)
$ws.Range("A5").Value = "greenJavaGreatSelenium"
$ws.Range("B5").Value = "allAboutLogic"
$ws.Range("A6").Value = "KWAdfd"
$ws.Range("B7").Value = "BASicgraet"

# Update the active selection on the sheet
$ws.Activate()
$ws.Range("B7").Select()
